# Auto-generated edit script
# Updates Betfair Back/Lay odds cells per the target diff for 2025-11-13 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 1.3
$ws.Range("J2").Value = 1.09
$ws.Range("K2").Value = 44
$ws.Range("W2").Value = 1.09

# Row 3
$ws.Range("F3").Value = 2.08
$ws.Range("G3").Value = 2.16
$ws.Range("I3").Value = 4.7
$ws.Range("J3").Value = 3.15
$ws.Range("L3").Value = 1.6
$ws.Range("N3").Value = 2.48
$ws.Range("O3").Value = 1.59
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.76
$ws.Range("V3").Value = 1.22
$ws.Range("AA3").Value = 130
$ws.Range("AB3").Value = 6.6
$ws.Range("AC3").Value = 7.8
$ws.Range("AE3").Value = 90
$ws.Range("AH3").Value = 34
$ws.Range("AI3").Value = 140
$ws.Range("AO3").Value = 150

# Row 4
$ws.Range("F4").Value = 1.71
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 1.5
$ws.Range("J4").Value = 1.09
$ws.Range("K4").Value = 32
$ws.Range("W4").Value = 1.05

# Row 5
$ws.Range("K5").Value = 3.95
$ws.Range("S5").Value = 4.5
$ws.Range("U5").Value = 1.71
$ws.Range("AF5").Value = 9.6
$ws.Range("AG5").Value = 10.5
$ws.Range("AH5").Value = 28

# Row 6
$ws.Range("F6").Value = 3.55
$ws.Range("G6").Value = 4.3
$ws.Range("H6").Value = 2.02
$ws.Range("I6").Value = 2.22
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 4.5
$ws.Range("P6").Value = 1.79
$ws.Range("Q6").Value = 1.88
$ws.Range("U6").Value = 1.96
$ws.Range("V6").Value = 1.83
$ws.Range("W6").Value = 1.32

# Row 7
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 1.29
$ws.Range("I7").Value = 1.36
$ws.Range("J7").Value = 6.2
$ws.Range("K7").Value = 7.4
$ws.Range("O7").Value = 1.15
$ws.Range("P7").Value = 2.76
$ws.Range("Q7").Value = 1.45
$ws.Range("T7").Value = 1.84
$ws.Range("U7").Value = 1.96
$ws.Range("W7").Value = 1.08
$ws.Range("Z7").Value = 10
$ws.Range("AB7").Value = 44
$ws.Range("AC7").Value = 15
$ws.Range("AH7").Value = 32
$ws.Range("AJ7").Value = 510
$ws.Range("AL7").Value = 150
$ws.Range("AM7").Value = 140

# Row 8
$ws.Range("F8").Value = 1.87
$ws.Range("G8").Value = 2.08
$ws.Range("H8").Value = 1.81
$ws.Range("I8").Value = 5.6
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 4.3
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 2.86
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 1.64
$ws.Range("Q8").Value = 1.94
$ws.Range("R8").Value = 1.24
$ws.Range("S8").Value = 3.4
$ws.Range("V8").Value = 1.22
$ws.Range("W8").Value = 1.92

# Row 9
$ws.Range("G9").Value = 2.44
$ws.Range("H9").Value = 3.1
$ws.Range("J9").Value = 3.3
$ws.Range("K9").Value = 3.75
$ws.Range("Q9").Value = 1.76
$ws.Range("S9").Value = 2.88
$ws.Range("T9").Value = 1.63
$ws.Range("V9").Value = 1.38
$ws.Range("W9").Value = 1.7
$ws.Range("AA9").Value = 65

# Row 10
$ws.Range("F10").Value = 2.82
$ws.Range("G10").Value = 3.1
$ws.Range("H10").Value = 2.44
$ws.Range("I10").Value = 2.68
$ws.Range("K10").Value = 3.7
$ws.Range("L10").Value = 1.34
$ws.Range("P10").Value = 1.88
$ws.Range("Q10").Value = 1.93
$ws.Range("T10").Value = 1.73
$ws.Range("U10").Value = 2.1
$ws.Range("W10").Value = 1.47
$ws.Range("X10").Value = 15
$ws.Range("Y10").Value = 11
$ws.Range("Z10").Value = 17.5
$ws.Range("AA10").Value = 1000
$ws.Range("AC10").Value = 8.4
$ws.Range("AD10").Value = 12.5
$ws.Range("AE10").Value = 1000
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 18.5
$ws.Range("AI10").Value = 980
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 120
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 25

# Row 11
$ws.Range("H11").Value = 2.94
$ws.Range("N11").Value = 2.78
$ws.Range("O11").Value = 1.43
$ws.Range("T11").Value = 1.89
$ws.Range("W11").Value = 1.51

# Row 12
$ws.Range("F12").Value = 4.1
$ws.Range("I12").Value = 2.16
$ws.Range("J12").Value = 3.2
$ws.Range("K12").Value = 3.65
$ws.Range("L12").Value = 1.42
$ws.Range("M12").Value = 1.09
$ws.Range("N12").Value = 2.82
$ws.Range("O12").Value = 1.43
$ws.Range("P12").Value = 1.66
$ws.Range("Q12").Value = 2.28
$ws.Range("S12").Value = 4.3
$ws.Range("T12").Value = 2
$ws.Range("U12").Value = 1.82
$ws.Range("V12").Value = 1.87
$ws.Range("W12").Value = 1.27
$ws.Range("X12").Value = 11
$ws.Range("Y12").Value = 8
$ws.Range("Z12").Value = 12.5
$ws.Range("AB12").Value = 13
$ws.Range("AE12").Value = 27
$ws.Range("AH12").Value = 23
$ws.Range("AL12").Value = 95

# Row 13
$ws.Range("F13").Value = 1.49
$ws.Range("G13").Value = 1.66
$ws.Range("H13").Value = 6.6
$ws.Range("I13").Value = 9.800000000000001
$ws.Range("J13").Value = 3.95
$ws.Range("K13").Value = 5.6
$ws.Range("L13").Value = 1.37
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 3.45
$ws.Range("O13").Value = 1.38
$ws.Range("P13").Value = 1.81
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 1.3
$ws.Range("S13").Value = 3.95
$ws.Range("T13").Value = 2.12
$ws.Range("U13").Value = 1.7
$ws.Range("V13").Value = 1.12
$ws.Range("W13").Value = 2.5
$ws.Range("AB13").Value = 8.6

# Row 14
$ws.Range("F14").Value = 2.26
$ws.Range("G14").Value = 2.46
$ws.Range("H14").Value = 3.55
$ws.Range("I14").Value = 4.1
$ws.Range("J14").Value = 3.15
$ws.Range("K14").Value = 3.4
$ws.Range("M14").Value = 1.1
$ws.Range("Q14").Value = 2.28
$ws.Range("R14").Value = 1.23
$ws.Range("T14").Value = 1.93
$ws.Range("U14").Value = 1.87
$ws.Range("V14").Value = 1.32
$ws.Range("W14").Value = 1.68
$ws.Range("X14").Value = 10.5
$ws.Range("Y14").Value = 12
$ws.Range("Z14").Value = 980
$ws.Range("AA14").Value = 100
$ws.Range("AB14").Value = 8.6
$ws.Range("AC14").Value = 7.8
$ws.Range("AD14").Value = 17
$ws.Range("AF14").Value = 14.5
$ws.Range("AG14").Value = 12.5
$ws.Range("AH14").Value = 21
$ws.Range("AJ14").Value = 980
$ws.Range("AK14").Value = 980
$ws.Range("AN14").Value = 980
$ws.Range("AO14").Value = 75

# Row 15
$ws.Range("F15").Value = 1.75
$ws.Range("G15").Value = 1.91
$ws.Range("H15").Value = 4.9
$ws.Range("I15").Value = 5.8
$ws.Range("J15").Value = 3.5
$ws.Range("M15").Value = 1.07
$ws.Range("P15").Value = 1.81
$ws.Range("Q15").Value = 1.98
$ws.Range("U15").Value = 1.9
$ws.Range("V15").Value = 1.2
$ws.Range("W15").Value = 2.08
$ws.Range("X15").Value = 15
$ws.Range("Y15").Value = 17.5
$ws.Range("Z15").Value = 980
$ws.Range("AA15").Value = 160
$ws.Range("AB15").Value = 8.6
$ws.Range("AC15").Value = 8.800000000000001
$ws.Range("AD15").Value = 22
$ws.Range("AH15").Value = 22
$ws.Range("AI15").Value = 960
$ws.Range("AJ15").Value = 21
$ws.Range("AL15").Value = 980
$ws.Range("AM15").Value = 150
$ws.Range("AN15").Value = 14.5

# Row 16
$ws.Range("G16").Value = 5
$ws.Range("K16").Value = 3.85
$ws.Range("N16").Value = 2.92
$ws.Range("R16").Value = 1.25
$ws.Range("W16").Value = 1.25
$ws.Range("Z16").Value = 13.5

# Row 17
$ws.Range("G17").Value = 3.85
$ws.Range("H17").Value = 2.2
$ws.Range("V17").Value = 1.77
$ws.Range("Y17").Value = 10
$ws.Range("AB17").Value = 16.5
$ws.Range("AK17").Value = 60
$ws.Range("AO17").Value = 19.5

# Row 18
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 1.45
$ws.Range("I18").Value = 1.56
$ws.Range("K18").Value = 5
$ws.Range("N18").Value = 3.25
$ws.Range("Q18").Value = 2
$ws.Range("S18").Value = 3.85
$ws.Range("T18").Value = 2.2
$ws.Range("U18").Value = 1.66
$ws.Range("V18").Value = 2.78
$ws.Range("W18").Value = 1.11

